$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has been recorded as no decision being made.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for the movie ""Barbie.""`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached on a movie for Friday.`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision to show ""Barbie"" on Friday has been recorded.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded as having no consensus on which movie to show on Friday.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to show ""Barbie.""`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie to be shown on Friday.`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision was made not to select a movie for Friday.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision process did not yield a consensus on which movie to show on Friday, resulting in no decision being made.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision on which movie to show on Friday.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded as ""no decision"" regarding the movie selection for Friday.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both movies.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired for the Friday showing.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday ended without a definitive choice.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision has been recorded as ""no decision"" regarding the movie to show on Friday.`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The movie ""Oppenheimer"" has been successfully selected for acquisition.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision to acquire the rights to both movies has been recorded successfully.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for ""Barbie"" will be acquired for the show on Friday.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision about Friday's movie can be made.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: There was no decision about the movie to play on Friday.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for the movie ""Barbie.""`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been made that no definitive choice about Friday's movie can be concluded at this time.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision to show a movie on Friday could not be made, so no action will be taken regarding acquiring movie rights.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been recorded as ""no decision"" regarding the movie selection for Friday.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been recorded and ""Barbie"" will be shown on Friday.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded as ""no decision."" If you need further assistance, feel free to ask!`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision-making process did not reach a conclusion about which movie to show on Friday, so the appropriate action is to indicate that no decision has been made.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision-making process resulted in no agreement on which movie to show on Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday is noted as undecided.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The function has been called successfully, and the decision outcome is recorded as no decision on which movie to show on Friday.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: No decision was made regarding which movie to show on Friday.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded as ""no decision"" regarding which movie to show on Friday.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday concluded without a selection.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: No movie was selected in this meeting.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision-making process has concluded without a choice of movie for Friday.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie"" to be shown on Friday.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday was not made, so no action has been taken.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision process concluded without a clear agreement on which movie to show on Friday. As a result, no movie rights will be acquired.`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision has been recorded as no selection for Friday's movie.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision regarding which movie to acquire has been recorded as no decision.`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The function for no decision has been executed. There was no agreement on what movie to show on Friday.`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both movies.`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision-making process ended without a definitive choice for a movie to be shown on Friday, so no movie has been acquired.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision has been recorded as ""no decision.""`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been recorded as that no movie will be shown on Friday.`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision to acquire the rights for ""Barbie"" has been made.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selection was made.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded.`n"
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be shown on Friday.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for ""Barbie.""`n"
$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for ""Barbie.""`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision to acquire the rights for ""Barbie"" has been recorded.`n"
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has resulted in no agreement.`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision concluded with no specific movie chosen for Friday.`n"
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday resulted in no agreement.`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C67").Value = "MSG: None`n`nMSG: The rights to both movies ""Oppenheimer"" and ""Barbie"" have been acquired.`n"
$ws.Range("C68").Value = "MSG: None`n`nMSG: The function for no decision has been called, indicating that no agreement was reached regarding which movie to show on Friday.`n"
$ws.Range("C69").Value = "MSG: None`n`nMSG: The decision has been recorded as a ""no decision,"" indicating that no agreement was reached on which movie to show on Friday.`n"
$ws.Range("C70").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both ""Barbie"" and ""Oppenheimer"" has been recorded successfully.`n"

$ws.Range("D7").Value = "Barbie_was_selected, "
$ws.Range("D14").Value = "both_movies, "
$ws.Range("D22").Value = "Barbie_was_selected, "
$ws.Range("D38").Value = "Barbie_was_selected, "
$ws.Range("D39").Value = "both_movies, "
$ws.Range("D43").Value = "both_movies, "
$ws.Range("D48").Value = "both_movies, "
